$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51, pushing the existing rows 51-55 down to 52-56.
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new weekly price entry.
$ws.Cells.Item(51, 1).Value = 11
$ws.Cells.Item(51, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(51, 3).Value = "Bíobío"
$ws.Cells.Item(51, 4).Value = 44783
$ws.Cells.Item(51, 5).Value = 8
$ws.Cells.Item(51, 6).Value = 100112031
$ws.Cells.Item(51, 7).Value = "Poroto verde"
$ws.Cells.Item(51, 8).Value = "Magnum"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 100
$ws.Cells.Item(51, 11).Value = 38000
$ws.Cells.Item(51, 12).Value = 40000
$ws.Cells.Item(51, 13).Value = 39000
$ws.Cells.Item(51, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(51, 15).Value = "Perú"
$ws.Cells.Item(51, 16).Value = 1560
$ws.Cells.Item(51, 17).Value = 25
$ws.Cells.Item(51, 18).Value = "Hortaliza"
